# Celerio-style regeneration of savedsearch.xlsx:
#  - rename the "accountId" property to "account" and route its value
#    through printer.print(...)
#  - add a matching descriptor row on the "Search" sheet

$wb = $excel.ActiveWorkbook

$listSheet   = $wb.Worksheets.Item("List")
$searchSheet = $wb.Worksheets.Item("Search")

$accountLabel   = '${msg.getProperty(''savedSearch_account'')}'
$accountFormula = '${printer.print(savedSearch.account)}'
$accountVar     = '${account}'

# "List" sheet: column E used to describe savedSearch.accountId, now describes
# savedSearch.account (printed through the printer helper).
$listSheet.Range("E1").Value = $accountLabel
$listSheet.Range("E2").Value = $accountFormula

# "Search" sheet: append the new descriptor row (row 7) for the account field.
$searchSheet.Range("A7").Value = $accountLabel
$searchSheet.Range("B7").Value = $accountVar
